$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34, shifting old rows 34-113 down to 35-114.
$ws.Rows("34:34").Insert()

# Fill the newly inserted (blank) row 34 with a copy of row 33's full content,
# since in the edited file row 34 duplicates row 33's data.
$ws.Range("A33:R33").Copy()
$ws.Range("A34").PasteSpecial(-4104)  # xlPasteAll

# Overwrite row 32 with the new weekly record's varying values.
$ws.Range("D32").Value = 45002
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 17000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17500
$ws.Range("P32").Value = 1346
